$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F ("dSF") values for specific rows per repull of data
$ws.Range("F4").Value = -1
$ws.Range("F9").Value = -7
$ws.Range("F11").Value = 2
$ws.Range("F16").Value = -4
